# daily auto push: 2026-03-01 03:22 UTC
# Insert a new data row for 2026/03/01 07:00 (day=日) at row 908,
# shifting all subsequent rows (old 908..949) down by one (new 909..950).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 908; existing rows 908-949 move to 909-950.
$ws.Rows.Item(908).Insert()

# Force column A to be stored as text so the date-like string "2026/03/01"
# isn't auto-converted into a date serial number by Excel.
$ws.Range("A908").NumberFormat = "@"

$ws.Cells.Item(908, 1).Value = "2026/03/01"
$ws.Cells.Item(908, 2).Value = "日"
$ws.Cells.Item(908, 3).Value = 7
$ws.Cells.Item(908, 4).Value = 201
